# Edit: insert two new quarterly columns (2018-12-31 and 2018-09-30) in front of the
# existing quarter data on the "GTS" sheet, shifting the prior quarters right by two
# columns (old D:K -> F:M), and populate the two new columns with the latest figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at D:E; existing data (D:K) shifts right to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# Give the new D:E columns the same number formatting as column F (the old column D,
# now shifted one place right) so dates/numbers render the same way as their neighbors.
# Done in the three contiguous data blocks (separated by genuinely blank spacer rows
# 36 and 78) so we don't manufacture phantom formatted-but-empty rows there.
$ws.Range("F5:F35").Copy()
$ws.Range("D5:E5").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New-quarter values: row number -> @(column D value, column E value)
$data = @{
    7 = @(43465, 43373)
    8 = @(705400, 770600)
    9 = @("NA", "NA")
    10 = @("NA", "NA")
    11 = @($null, $null)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @("NA", "NA")
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(713800, 789600)
    18 = @(-8400, -19000)
    19 = @($null, $null)
    20 = @(0, 0)
    21 = @(-4800, -16100)
    22 = @(1400, 2000)
    23 = @(-9800, -21000)
    24 = @(1100, -3400)
    25 = @(0, 0)
    26 = @(-10900, -17600)
    27 = @(-10900, -17600)
    28 = @(0, 0)
    29 = @(0, 0)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(0, 0)
    33 = @(-10900, -17600)
    34 = @(0, 0)
    35 = @(-10900, -17600)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(117500, 107100)
    42 = @(0, 0)
    43 = @(628400, 632900)
    44 = @(0, 0)
    45 = @(0, 0)
    46 = @(0, 0)
    47 = @(1564500, 1615700)
    48 = @(81900, 78400)
    49 = @(25400, 25400)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(79000, 83600)
    53 = @(0, 0)
    54 = @(2760200, 2818400)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(0, 0)
    58 = @(0, 0)
    59 = @(532300, 499400)
    60 = @(0, 0)
    61 = @(28900, 29700)
    62 = @(34500, 34100)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(1938300, 1997900)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(762000, 772900)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(822000, 820500)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(-10900, -17600)
    82 = @($null, $null)
    83 = @(3600, 2900)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(11400, -134700)
    90 = @($null, $null)
    91 = @(-7500, -3200)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(32700, -14500)
    95 = @($null, $null)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-33600, 300)
    101 = @(0, 0)
    102 = @(10500, -148900)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    if ($null -ne $vals[0]) {
        $ws.Cells.Item($r, 4).Value = $vals[0]
    }
    if ($null -ne $vals[1]) {
        $ws.Cells.Item($r, 5).Value = $vals[1]
    }
}
